$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fix tiny floating point precision differences in row 13 ---
$ws.Range("E13").Value = 0.9903279118225773
$ws.Range("I13").Value = 0.9892063469786482
$ws.Range("N13").Value = 0.990366156139844
$ws.Range("O13").Value = 0.9896128025487396

# --- Add new row 16 (HexGrid-60degTilt5degRes) ---
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 1.204979639723102
$ws.Range("D16").Value = 1.130889471180497
$ws.Range("E16").Value = 0.9868289383519783
$ws.Range("F16").Value = 0.9858307263132607
$ws.Range("G16").Value = 1.204979639723102
$ws.Range("H16").Value = 1.130889471180497
$ws.Range("I16").Value = 0.9600196664009265
$ws.Range("J16").Value = 0.9101211500444905
$ws.Range("K16").Value = 1.037813908353118
$ws.Range("L16").Value = 1.012825126442885
$ws.Range("M16").Value = 1.203803169134867
$ws.Range("N16").Value = 1.058859204766238
$ws.Range("O16").Value = 1.07713219389221
$ws.Range("P16").Value = 1.028663578351282

# Copy the formatting from row 15 (A15) onto A16, to match the bordered /
# bold / centered style used for every cell in column A.
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
